$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to retain a numeric-looking string as text,
# matching the original inline-string cell formatting, without leaving
# a lingering custom style on the cell (NumberFormat reset via Style).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- Row 41 / 42: Stacks and Monero swap position/rank ---
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D41") "176.37"
$ws.Range("E41").Value = "  -0.32%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D42") "2.21"
$ws.Range("E42").Value = "  +3.06%  "

# --- Price (column D) updates ---
$ws.Range("D2").Value = "67.813.80"
$ws.Range("D3").Value = "3.544.12"
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("D5") "612.77"
Set-TextValue $ws.Range("D6") "152.53"
$ws.Range("D7").Value = "3.542.51"
$ws.Range("D14").Value = "4.148.19"
Set-TextValue $ws.Range("D15") "31.97"
$ws.Range("D16").Value = "3.547.95"
$ws.Range("D17").Value = "67.605.64"
Set-TextValue $ws.Range("D20") "15.22"
Set-TextValue $ws.Range("D21") "9.68"
Set-TextValue $ws.Range("D22") "446.55"
Set-TextValue $ws.Range("D23") "0.623"
Set-TextValue $ws.Range("D24") "77.01"
$ws.Range("D26").Value = "3.685.40"
Set-TextValue $ws.Range("D28") "10.19"
Set-TextValue $ws.Range("D29") "8.62"
Set-TextValue $ws.Range("D32") "0.168"
Set-TextValue $ws.Range("D34") "25.76"
Set-TextValue $ws.Range("D35") "6.21"
$ws.Range("D36").Value = "3.531.26"
Set-TextValue $ws.Range("D38") "8.04"
Set-TextValue $ws.Range("D40") "1.00"
Set-TextValue $ws.Range("D43") "0.0894"
Set-TextValue $ws.Range("D46") "28.77"
Set-TextValue $ws.Range("D47") "45.57"
Set-TextValue $ws.Range("D48") "2.67"
Set-TextValue $ws.Range("D50") "7.60"

# --- Volume(1h) percentage (column E) updates ---
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("E11").Value = "  +3.38%  "
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("E20").Value = "  -1.71%  "
$ws.Range("E21").Value = "  +3.15%  "
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("E23").Value = "  -2.47%  "
$ws.Range("E24").Value = "  -2.43%  "
$ws.Range("E25").Value = "  +5.28%  "
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("E29").Value = "  +3.24%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  -4.13%  "
$ws.Range("E32").Value = "  +6.25%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("E34").Value = "  -0.71%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("E37").Value = "  -2.61%  "
$ws.Range("E38").Value = "  -0.90%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("E43").Value = "  +2.04%  "
$ws.Range("E44").Value = "  -3.69%  "
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("E49").Value = "  +3.38%  "
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("E51").Value = "  +0.25%  "
